# Update schedule and render
# - Shorten the "Differential Gene Expression" lecture block on day 2 (row 24)
#   so it ends at 11:30 instead of 12:00, and insert a new row for the
#   "Gene set analysis" lecture (11:30-12:00) right after it. This pushes
#   every subsequent row of the schedule down by one.
# - Update the "schedule" defined name so it still spans the whole table.
# - Leave the selection where the author last left it (F47).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shrink the Differential Gene Expression lecture to end at 11:30 (was 12:00)
$ws.Range("D24").Value = 0.47916666666666669

# Insert a new row for the Gene set analysis lecture (11:30-12:00)
$ws.Rows.Item(25).Insert()

$ws.Range("C25").Value = 0.47916666666666669
$ws.Range("D25").Value = 0.5
$ws.Range("E25").Value = "Lecture: Gene set analysis"
$ws.Range("F25").Value = "Jennifer Fransson"
$ws.Range("H25").Value = "lectures/gsa/index.html"

# The old "Gene set analysis" row (now shifted down to row 28 by the
# insert above) is repurposed into an extra "Lab: Differential
# expression" slot, pointing at the lab instructions page instead of
# the lecture slides.
$ws.Range("E28").Value = "Lab: Differential expression"
$ws.Range("H28").ClearContents()
$ws.Range("I28").Value = "home_contents.html"

# Extend the named range "schedule" to cover the newly added row
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!schedule") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$L`$42"
    }
}

# Restore the active cell selection as last left by the author
$ws.Range("F47").Select()
